$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.804.70"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "3.417.92"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "570.52"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").Value = "158.24"
$ws.Range("E6").Value = "  +1.02%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.420.34"
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  -9.27%  "
$ws.Range("D10").Value = "7.26"
$ws.Range("E10").Value = "  +1.84%  "
$ws.Range("E11").Value = "  -2.95%  "
$ws.Range("E12").Value = "  -4.12%  "
$ws.Range("D13").Value = "4.005.97"
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("D15").Value = "27.10"
$ws.Range("E15").Value = "  -2.08%  "
$ws.Range("E16").Value = "  -8.00%  "
$ws.Range("D17").Value = "63.878.44"
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("D18").Value = "3.408.76"
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("E19").Value = "  -3.45%  "
$ws.Range("E20").Value = "  -1.74%  "
$ws.Range("D21").Value = "381.70"
$ws.Range("E22").Value = "  -2.51%  "
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").Value = "71.27"
$ws.Range("E24").Value = "  -0.68%  "
$ws.Range("E25").Value = "  -5.85%  "
$ws.Range("E26").Value = "  -2.70%  "
$ws.Range("D27").Value = "9.67"
$ws.Range("E27").Value = "  -5.99%  "
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  -1.35%  "
$ws.Range("E31").Value = "  -5.71%  "
$ws.Range("D32").Value = "1.98"
$ws.Range("E32").Value = "  -0.92%  "
$ws.Range("D34").Value = "22.86"
$ws.Range("E34").Value = "  -0.84%  "
$ws.Range("E35").Value = "  -2.55%  "
$ws.Range("E36").Value = "  -6.11%  "
$ws.Range("D37").Value = "160.76"
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("E38").Value = "  +9.71%  "
$ws.Range("E39").Value = "  -3.70%  "
$ws.Range("D40").Value = "2.812.17"
$ws.Range("E40").Value = "  -2.17%  "
$ws.Range("D41").Value = "26.06"
$ws.Range("E41").Value = "  -0.94%  "
$ws.Range("E42").Value = "  -5.26%  "
$ws.Range("D43").Value = "43.07"
$ws.Range("E43").Value = "  +0.50%  "
$ws.Range("D44").Value = "26.33"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("E45").Value = "  -8.02%  "
$ws.Range("E46").Value = "  -5.92%  "
$ws.Range("E47").Value = "  -3.64%  "
$ws.Range("D48").Value = "2.40"
$ws.Range("E48").Value = "  +9.67%  "
$ws.Range("D49").Value = "333.18"
$ws.Range("E49").Value = "  +3.68%  "
$ws.Range("E50").Value = "  -3.89%  "
$ws.Range("E51").Value = "  -5.35%  "
